# Auto-generated: updates currentAveragePrice / LevePrice / LeveProfit
# columns (H-N) for specific leve rows across multiple sheets,
# reflecting refreshed market-board pricing data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 8: On the Drip / Eye Drops
$ws.Range("H8").Value = 71.875
$ws.Range("I8").Value = 71.875
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 215.625
$ws.Range("L8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -76.625

# Row 19: Unbreak My Heart / Roof Tile
$ws.Range("H19").Value = 201.41379
$ws.Range("I19").Value = 239.08333
$ws.Range("J19").Value = 174.82353
$ws.Range("K19").Value = 239.08333
$ws.Range("L19").Value = 174.82353
$ws.Range("M19").Value = -64.08332999999999
$ws.Range("N19").Value = -524.82353

# Row 69: Steeling the Knife, Steeling the Mind / Grade 1 Mind Dissolvent
$ws.Range("H69").Value = 4521.8696
$ws.Range("J69").Value = 4928.7856
$ws.Range("L69").Value = 14786.3568
$ws.Range("N69").Value = -16534.3568

# Row 72: Surgical Substitution (L) / Grade 1 Mind Dissolvent
$ws.Range("H72").Value = 4521.8696
$ws.Range("J72").Value = 4928.7856
$ws.Range("L72").Value = 44359.0704
$ws.Range("N72").Value = -53095.0704

# Row 74: Adhesive of Antipathy / Wing Glue
$ws.Range("H74").Value = 5950.3076
$ws.Range("I74").Value = 3300
$ws.Range("J74").Value = 8222
$ws.Range("K74").Value = 3300
$ws.Range("L74").Value = 8222
$ws.Range("M74").Value = -2364
$ws.Range("N74").Value = -10094

# Row 77: It's Gonna Grow Back (L) / Wing Glue
$ws.Range("H77").Value = 5950.3076
$ws.Range("I77").Value = 3300
$ws.Range("J77").Value = 8222
$ws.Range("K77").Value = 16500
$ws.Range("L77").Value = 41110
$ws.Range("M77").Value = -11820
$ws.Range("N77").Value = -50470

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 1597.8438
$ws.Range("I137").Value = 1107.8125
$ws.Range("J137").Value = 2087.875
$ws.Range("K137").Value = 3323.4375
$ws.Range("L137").Value = 6263.625
$ws.Range("M137").Value = -773.4375
$ws.Range("N137").Value = -11363.625

$ws = $wb.Worksheets.Item("ARM")
# Row 6: Don't Hit Me One More Time / Bronze Hoplon
$ws.Range("H6").Value = 641801.2
$ws.Range("I6").Value = 1000
$ws.Range("J6").Value = 802001.5
$ws.Range("K6").Value = 1000
$ws.Range("L6").Value = 802001.5
$ws.Range("M6").Value = -827
$ws.Range("N6").Value = -802347.5

# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 2498
$ws.Range("I61").Value = 2249.4783
$ws.Range("J61").Value = 3212.5
$ws.Range("K61").Value = 2249.4783
$ws.Range("L61").Value = 3212.5
$ws.Range("M61").Value = -2037.4783
$ws.Range("N61").Value = -3636.5

# Row 63: Rivets Run through It / Mythrite Rivets
$ws.Range("H63").Value = 100001870
$ws.Range("I63").Value = 142858830
$ws.Range("J63").Value = 2300
$ws.Range("K63").Value = 142858830
$ws.Range("L63").Value = 2300
$ws.Range("M63").Value = -142858144
$ws.Range("N63").Value = -3672

# Row 66: A Riveting Revival (L) / Mythrite Rivets
$ws.Range("H66").Value = 100001870
$ws.Range("I66").Value = 142858830
$ws.Range("J66").Value = 2300
$ws.Range("K66").Value = 714294150
$ws.Range("L66").Value = 11500
$ws.Range("M66").Value = -714290718
$ws.Range("N66").Value = -18364

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 2177
$ws.Range("I132").Value = 1936.2894
$ws.Range("J132").Value = 3701.5
$ws.Range("K132").Value = 5808.8682
$ws.Range("L132").Value = 11104.5
$ws.Range("M132").Value = -3278.8682
$ws.Range("N132").Value = -16164.5

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 2498
$ws.Range("I136").Value = 2249.4783
$ws.Range("J136").Value = 3212.5
$ws.Range("K136").Value = 6748.4349
$ws.Range("L136").Value = 9637.5
$ws.Range("M136").Value = -4198.4349
$ws.Range("N136").Value = -14737.5

$ws = $wb.Worksheets.Item("BSM")
# Row 12: A Hit Job / Bronze Chaser Hammer
$ws.Range("H12").Value = 722.5
$ws.Range("I12").Value = 722.5
$ws.Range("K12").Value = 722.5
$ws.Range("M12").Value = -554.5

$ws = $wb.Worksheets.Item("CUL")
# Row 37: I Love Lamprey / Eel Pie
$ws.Range("H37").Value = 42111.11
$ws.Range("J37").Value = 42111.11
$ws.Range("L37").Value = 126333.33
$ws.Range("N37").Value = -126557.33

# Row 68: Such a Butter Face / Fermented Butter
$ws.Range("H68").Value = 3328.9565
$ws.Range("I68").Value = 4535.1113
$ws.Range("J68").Value = 1614.9474
$ws.Range("K68").Value = 13605.3339
$ws.Range("L68").Value = 4844.8422
$ws.Range("M68").Value = -12794.3339
$ws.Range("N68").Value = -6466.8422

# Row 70: Persona non Gratin / Dhalmel Gratin
$ws.Range("H70").Value = 3137.353
$ws.Range("I70").Value = 1537.2222
$ws.Range("J70").Value = 4937.5
$ws.Range("K70").Value = 4611.6666
$ws.Range("L70").Value = 14812.5
$ws.Range("M70").Value = -4296.6666
$ws.Range("N70").Value = -15442.5

# Row 71: No Margarine of Error (L) / Fermented Butter
$ws.Range("H71").Value = 3328.9565
$ws.Range("I71").Value = 4535.1113
$ws.Range("J71").Value = 1614.9474
$ws.Range("K71").Value = 40816.00169999999
$ws.Range("L71").Value = 14534.5266
$ws.Range("M71").Value = -36760.00169999999
$ws.Range("N71").Value = -22646.5266

# Row 73: Recipe for Disaster (L) / Dhalmel Gratin
$ws.Range("H73").Value = 3137.353
$ws.Range("I73").Value = 1537.2222
$ws.Range("J73").Value = 4937.5
$ws.Range("K73").Value = 4611.6666
$ws.Range("L73").Value = 14812.5
$ws.Range("M73").Value = -3519.6666
$ws.Range("N73").Value = -16996.5

# Row 107: Slippery Service / Frantoio Oil
$ws.Range("H107").Value = 622.39703
$ws.Range("I107").Value = 199.5
$ws.Range("J107").Value = 824.65216
$ws.Range("K107").Value = 598.5
$ws.Range("L107").Value = 2473.95648
$ws.Range("M107").Value = 1321.5
$ws.Range("N107").Value = -6313.95648

# Row 130: Blast from the Pasta / The Noodles of Elpis
$ws.Range("H130").Value = 5690
$ws.Range("I130").Value = 2980
$ws.Range("J130").Value = 5936.364
$ws.Range("K130").Value = 8940
$ws.Range("L130").Value = 17809.092
$ws.Range("M130").Value = -3920
$ws.Range("N130").Value = -27849.092

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 18334412
$ws.Range("I131").Value = 6250598
$ws.Range("J131").Value = 22728526
$ws.Range("K131").Value = 18751794
$ws.Range("L131").Value = 68185578
$ws.Range("M131").Value = -18746754
$ws.Range("N131").Value = -68195658

# Row 133: Friends Are Food / Boiled Alpaca Steak
$ws.Range("H133").Value = 53595.047
$ws.Range("I133").Value = 172348.5
$ws.Range("J133").Value = 9062.5
$ws.Range("K133").Value = 517045.5
$ws.Range("L133").Value = 27187.5
$ws.Range("M133").Value = -511985.5
$ws.Range("N133").Value = -37307.5

# Row 134: Don't Knock It Till You've Tried It / Mezcal-marinated Swampmonk
$ws.Range("H134").Value = 10421.813
$ws.Range("I134").Value = 12113.9
$ws.Range("J134").Value = 9909.061
$ws.Range("K134").Value = 36341.7
$ws.Range("L134").Value = 29727.183
$ws.Range("M134").Value = -31271.7
$ws.Range("N134").Value = -39867.183

# Row 137: Creative Chocolate / Gateau au Chocolat
$ws.Range("H137").Value = 37052532
$ws.Range("I137").Value = 2133.6365
$ws.Range("J137").Value = 62524680
$ws.Range("K137").Value = 6400.9095
$ws.Range("L137").Value = 187574040
$ws.Range("M137").Value = -1300.9095
$ws.Range("N137").Value = -187584240

# Row 139: Najoothie / Wild Banana Blend
$ws.Range("H139").Value = 54114.24
$ws.Range("I139").Value = 92599.914
$ws.Range("K139").Value = 277799.742
$ws.Range("M139").Value = -272659.742

$ws = $wb.Worksheets.Item("GSM")
# Row 31: One and Only / Staghorn Staff
$ws.Range("H31").Value = 2221.8333
$ws.Range("I31").Value = 2221.8333
$ws.Range("K31").Value = 2221.8333
$ws.Range("M31").Value = -1929.8333

# Row 37: Dancing with the Stars / Toothed Staghorn Staff
$ws.Range("H37").Value = 2221.8333
$ws.Range("I37").Value = 2221.8333
$ws.Range("K37").Value = 2221.8333
$ws.Range("M37").Value = -1944.8333

# Row 52: It's My Business to Know Things / Red Coral Armillae
$ws.Range("H52").Value = 24350
$ws.Range("J52").Value = 24350
$ws.Range("L52").Value = 24350
$ws.Range("N52").Value = -24868

# Row 80: Needs More Prayerbell / Hardsilver Ingot
$ws.Range("H80").Value = 2646.3333
$ws.Range("I80").Value = 2670.5881
$ws.Range("K80").Value = 2670.5881
$ws.Range("M80").Value = -1672.5881

# Row 82: Appeasing the Astromancer / Hardsilver Planisphere
$ws.Range("H82").Value = 42330
$ws.Range("J82").Value = 42330
$ws.Range("L82").Value = 42330
$ws.Range("N82").Value = -43096

# Row 83: With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Range("H83").Value = 2646.3333
$ws.Range("I83").Value = 2670.5881
$ws.Range("K83").Value = 13352.9405
$ws.Range("M83").Value = -8360.940500000001

# Row 85: Silver Bar of Upcycling (L) / Hardsilver Planisphere
$ws.Range("H85").Value = 42330
$ws.Range("J85").Value = 42330
$ws.Range("L85").Value = 42330
$ws.Range("N85").Value = -44982

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs / Aldgoat Leather
$ws.Range("H22").Value = 2780376
$ws.Range("J22").Value = 2911.1765
$ws.Range("L22").Value = 2911.1765
$ws.Range("N22").Value = -3501.1765

# Row 27: Fire and Hide / Aldgoat Leather
$ws.Range("H27").Value = 2780376
$ws.Range("J27").Value = 2911.1765
$ws.Range("L27").Value = 2911.1765
$ws.Range("N27").Value = -3125.1765

# Row 81: I Need Your Glove Tonight / Dragonskin Gloves of Healing
$ws.Range("H81").Value = 33500
$ws.Range("J81").Value = 33500
$ws.Range("L81").Value = 33500
$ws.Range("N81").Value = -35496

# Row 84: Halonic Drake Handlers (L) / Dragonskin Gloves of Healing
$ws.Range("H84").Value = 33500
$ws.Range("J84").Value = 33500
$ws.Range("L84").Value = 100500
$ws.Range("N84").Value = -110484

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 21672402
$ws.Range("I132").Value = 25496002
$ws.Range("K132").Value = 76488006
$ws.Range("M132").Value = -76485476

$ws = $wb.Worksheets.Item("WVR")
# Row 11: Wiggle Room / Hempen Shepherd's Tunic
$ws.Range("H11").Value = 680000
$ws.Range("J11").Value = 680000
$ws.Range("L11").Value = 680000
$ws.Range("N11").Value = -680284

# Row 57: Glad As a Hatter / Felt Hat
$ws.Range("H57").Value = 50000
$ws.Range("I57").Value = 50000
$ws.Range("K57").Value = 50000
$ws.Range("M57").Value = -49246

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 1778.5862
$ws.Range("I132").Value = 1864.6154
$ws.Range("J132").Value = 1033
$ws.Range("K132").Value = 5593.8462
$ws.Range("L132").Value = 3099
$ws.Range("M132").Value = -3063.8462
$ws.Range("N132").Value = -8159

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 1031.3334
$ws.Range("I136").Value = 921.05884
$ws.Range("K136").Value = 2763.17652
$ws.Range("M136").Value = -213.17652

Write-Host "Sheets updated via scheduled runner."
